$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column AN into the new column AO
$ws.Range("AN1:AN54").Copy() | Out-Null
$ws.Range("AO1:AO54").PasteSpecial(-4122) | Out-Null

# New header for 4/22/20
$ws.Range("AO1").Value = " 4/22/20"

# New daily death counts per state
$values = @{
    2 = 201
    3 = 9
    4 = 229
    5 = 44
    6 = 1437
    7 = 508
    8 = 1544
    9 = 89
    10 = 127
    11 = 927
    12 = 846
    13 = 5
    14 = 12
    15 = 54
    16 = 1565
    17 = 661
    18 = 90
    19 = 111
    20 = 185
    21 = 1473
    22 = 39
    23 = 698
    24 = 2182
    25 = 2813
    26 = 179
    27 = 193
    28 = 232
    29 = 14
    30 = 42
    31 = 172
    32 = 48
    33 = 5063
    34 = 71
    35 = 20354
    36 = 265
    37 = 14
    38 = 610
    39 = 170
    40 = 78
    41 = 1713
    42 = 67
    43 = 181
    44 = 140
    45 = 9
    46 = 166
    47 = 550
    48 = 34
    49 = 40
    50 = 349
    51 = 692
    52 = 29
    53 = 246
    54 = 6
}
foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 41).Value = $values[$row]
}

$ws.Range("AO2").Select() | Out-Null

Write-Output "Added column AO ( 4/22/20) with updated state death totals."
